$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D..K shift to F..M)
$ws.Range("D:E").Insert()

# Copy number formats from column F (the old column D) across all used rows
# onto the two newly inserted columns D and E so dates/numbers format correctly.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)

$data = @{
    7 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    8 = @(1151100, 1276500, 1259200, 1215300, 1187200, 1221800, 1216600, 1220200, 1249200, 1244500)
    9 = @(1002000, 1052200, 1042800, 1064800, 983100, 1008100, 998200, 978000, 963500, 973800)
    10 = @(149100, 224300, 216300, 150400, 204100, 213700, 218400, 242200, 285600, 270700)
    11 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    12 = @(121400, 107900, 98100, 94800, 100100, 108400, 105200, 130100, 124700, 110800)
    13 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    14 = @(1500, 0, 0, 0, 8200, 9700, 3400, 9400, 34900, 17000)
    15 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    16 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    17 = @(1171700, 1197600, 1156100, 1190300, 1133800, 1178700, 1165900, 1184800, 1209800, 1198300)
    18 = @(-20600, 78900, 103100, 24900, 53400, 43100, 50700, 35400, 39300, 46200)
    19 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    20 = @(-40800, -29500, -12500, 56900, 25000, 70200, 36700, 16600, 11300, 29700)
    21 = @(340900, 469900, 524000, 512500, 507000, 550400, 511800, 485400, 486800, 496600)
    22 = @(22400, 22600, 23000, 21700, 21700, 20400, 18800, 17200, 15800, 12300)
    23 = @(-83800, 26900, 67600, 60200, 56700, 92900, 68600, 34800, 34800, 63600)
    24 = @(13400, 20500, -10700, -38000, 18000, 13000, 20700, -14000, 20100, 6300)
    25 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    26 = @(-97200, 6400, 78300, 98200, 38700, 79900, 47900, 48800, 14600, 57300)
    27 = @(-55300, 55800, 118600, 110200, 57400, 112500, 68000, 74500, 83100, 97000)
    28 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    29 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    30 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    31 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    32 = @(40800, 29500, 12500, -56900, -25000, -70200, -36700, -16600, -11300, -29700)
    33 = @(-55300, 55800, 118600, 110200, 57400, 112500, 68000, 74500, 83100, 97000)
    34 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    35 = @(-55300, 55800, 118600, 110200, 57400, 112500, 68000, 74500, 83100, 97000)
    38 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    39 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    40 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    41 = @(2711500, 2642100, 2437000, 2500200, 2647100, 2266700, 2208200, 1983100, 1877700, 1802200)
    42 = @(17100, 13700, 15200, 17000, 23200, 27400, 25900, 28300, 23300, 24500)
    43 = @(800400, 883300, 935000, 864000, 734300, 785600, 748200, 705300, 782700, 761800)
    44 = @(590000, 570000, 572200, 555400, 591700, 554300, 527800, 530300, 554300, 559800)
    45 = @(457100, 433300, 483400, 484800, 513800, 407900, 381500, 390200, 364400, 340000)
    46 = @(4576100, 4542400, 4442800, 4421400, 4510200, 4041800, 3891700, 3637200, 3602400, 3488300)
    47 = @(1085900, 1175200, 1209400, 1168800, 1102700, 1125500, 1163800, 1138100, 1133700, 1232600)
    48 = @(5623400, 5801300, 6147600, 6409700, 6677400, 6764100, 6921700, 7101400, 7375100, 7060600)
    49 = @(97000, 103900, 116100, 119500, 122800, 125500, 123900, 120300, 133300, 138400)
    50 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    51 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    52 = @(434500, 455400, 428700, 376600, 359800, 394100, 403200, 374500, 364200, 307400)
    53 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    54 = @(11816800, 12078200, 12344700, 12496000, 12772800, 12451000, 12504300, 12371600, 12608800, 12227400)
    55 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    56 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    57 = @(220400, 223200, 240100, 226800, 211800, 214200, 215400, 210000, 223500, 219900)
    58 = @(590700, 658800, 558500, 1296300, 1711600, 1579800, 1561900, 1185000, 1012600, 1534200)
    59 = @(806100, 785000, 1119200, 828900, 930700, 875200, 1064900, 785700, 1111100, 967000)
    60 = @(1617200, 1667000, 1917700, 2352000, 2854100, 2669200, 2842200, 2180600, 2347200, 2721000)
    61 = @(2174200, 2186300, 2196200, 1705200, 1728100, 1696900, 1660600, 2032000, 1980400, 1294100)
    62 = @(1331600, 1335900, 1399200, 1431400, 1253700, 1135100, 1091800, 1092900, 1148100, 1070700)
    63 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    64 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    65 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    66 = @(5138100, 5205200, 5533000, 5495200, 5866800, 5535000, 5624800, 5336500, 5546200, 5138300)
    67 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    68 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    69 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    70 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    71 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    72 = @(1996100, 2097900, 2074700, 2248700, 1557800, 1529000, 1457800, 1603500, 1554000, 1561100)
    73 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    74 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    75 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    76 = @(6678700, 6873000, 6811700, 7000800, 6905900, 6916000, 6879500, 7035100, 7062700, 7089100)
    77 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    80 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    81 = @(-55300, 55800, 118600, 110200, 57400, 112500, 68000, 74500, 83100, 97000)
    82 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    83 = @(402300, 420500, 433400, 430700, 428700, 437100, 424300, 433500, 436200, 420700)
    84 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    85 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    86 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    87 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    88 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    89 = @(392900, 511200, 462300, 284400, 418100, 371000, 480200, 434100, 507000, 452900)
    90 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    91 = @(-141300, -181900, -126400, -185300, -355300, -236200, -270100, -588000, -724600, -666200)
    92 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    93 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    94 = @(-149500, -177500, -196900, 21600, -195500, -161500, -262700, -531300, -597400, -502200)
    95 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    96 = @(0, 0, 0, 0, 300, 0, 0, 0, 0, -225200)
    97 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    98 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    99 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    100 = @(-185500, -105400, -349900, -444500, 169000, -149200, 22900, 255800, 137600, 289700)
    101 = @(11500, -23200, 21300, -8400, -11100, -1800, -3100, -53200, 28300, -49900)
    102 = @(69400, 205100, -63200, -146900, 380400, 58500, 237300, 105400, 75400, 190500)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        if ($null -ne $vals[$i]) {
            $ws.Cells.Item([int]$r, 4 + $i).Value = $vals[$i]
        }
    }
}
